# Auto-generated edit script applying the Sophia_Profits.xlsx diff
# (workbook sheet names: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3702.8572
$ws.Range("I17").Value = 2100
$ws.Range("J17").Value = 3970
$ws.Range("K17").Value = 6300
$ws.Range("L17").Value = 11910
$ws.Range("M17").Value = -6132
$ws.Range("N17").Value = -12246
$ws.Range("H40").Value = 4524.875
$ws.Range("I40").Value = 2699.8333
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 2699.8333
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -2524.8333
$ws.Range("N40").Value = -10350
$ws.Range("H43").Value = 1237.25
$ws.Range("J43").Value = 933.3333
$ws.Range("L43").Value = 933.3333
$ws.Range("N43").Value = -1071.3333
$ws.Range("H53").Value = 227
$ws.Range("I53").Value = 159
$ws.Range("K53").Value = 159
$ws.Range("M53").Value = 478
$ws.Range("H103").Value = 1761.2727
$ws.Range("I103").Value = 1368.75
$ws.Range("J103").Value = 1985.5714
$ws.Range("K103").Value = 4106.25
$ws.Range("L103").Value = 5956.7142
$ws.Range("M103").Value = -3520.25
$ws.Range("N103").Value = -7128.7142
$ws.Range("H138").Value = 4412.439
$ws.Range("I138").Value = 2305
$ws.Range("J138").Value = 5005.1562
$ws.Range("K138").Value = 6915
$ws.Range("L138").Value = 15015.4686
$ws.Range("M138").Value = -1775
$ws.Range("N138").Value = -25295.4686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8419.789000000001
$ws.Range("I32").Value = 8063.9443
$ws.Range("J32").Value = 14825
$ws.Range("K32").Value = 8063.9443
$ws.Range("L32").Value = 14825
$ws.Range("M32").Value = -7776.9443
$ws.Range("N32").Value = -15399
$ws.Range("H45").Value = 2475.6155
$ws.Range("I45").Value = 2518.3
$ws.Range("J45").Value = 2333.3333
$ws.Range("K45").Value = 2518.3
$ws.Range("L45").Value = 2333.3333
$ws.Range("M45").Value = -2141.3
$ws.Range("N45").Value = -3087.3333
$ws.Range("H61").Value = 2414.0588
$ws.Range("I61").Value = 2197.375
$ws.Range("J61").Value = 5881
$ws.Range("K61").Value = 2197.375
$ws.Range("L61").Value = 5881
$ws.Range("M61").Value = -1985.375
$ws.Range("N61").Value = -6305
$ws.Range("H135").Value = 53300
$ws.Range("J135").Value = 53300
$ws.Range("L135").Value = 53300
$ws.Range("N135").Value = -63440
$ws.Range("H136").Value = 2414.0588
$ws.Range("I136").Value = 2197.375
$ws.Range("J136").Value = 5881
$ws.Range("K136").Value = 6592.125
$ws.Range("L136").Value = 17643
$ws.Range("M136").Value = -4042.125
$ws.Range("N136").Value = -22743

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5122.273
$ws.Range("I20").Value = 1838.4445
$ws.Range("K20").Value = 1838.4445
$ws.Range("M20").Value = -1591.4445
$ws.Range("H22").Value = 87.5
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 73

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 492.85715
$ws.Range("I22").Value = 410.16666
$ws.Range("K22").Value = 410.16666
$ws.Range("M22").Value = -60.16665999999998
$ws.Range("H62").Value = 7799.2
$ws.Range("I62").Value = 7856
$ws.Range("K62").Value = 7856
$ws.Range("M62").Value = -7232
$ws.Range("H65").Value = 7799.2
$ws.Range("I65").Value = 7856
$ws.Range("K65").Value = 39280
$ws.Range("M65").Value = -36160
$ws.Range("H95").Value = 30099.8
$ws.Range("J95").Value = 30099.8
$ws.Range("L95").Value = 30099.8
$ws.Range("N95").Value = -35591.8
$ws.Range("H141").Value = 225994.14
$ws.Range("J141").Value = 225994.14
$ws.Range("L141").Value = 225994.14
$ws.Range("N141").Value = -236354.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.285713
$ws.Range("I2").Value = 25.347826
$ws.Range("J2").Value = 64.2
$ws.Range("K2").Value = 152.086956
$ws.Range("L2").Value = 385.2
$ws.Range("M2").Value = -39.08695600000001
$ws.Range("N2").Value = -611.2
$ws.Range("H8").Value = 6005153
$ws.Range("I8").Value = 6005153
$ws.Range("K8").Value = 18015459
$ws.Range("M8").Value = -18015320
$ws.Range("H14").Value = 1229.6666
$ws.Range("I14").Value = 1229.6666
$ws.Range("K14").Value = 3688.9998
$ws.Range("M14").Value = -3515.9998
$ws.Range("H80").Value = 12074.5
$ws.Range("I80").Value = 6497.5
$ws.Range("J80").Value = 13189.9
$ws.Range("K80").Value = 19492.5
$ws.Range("L80").Value = 39569.7
$ws.Range("M80").Value = -18556.5
$ws.Range("N80").Value = -41441.7
$ws.Range("H83").Value = 12074.5
$ws.Range("I83").Value = 6497.5
$ws.Range("J83").Value = 13189.9
$ws.Range("K83").Value = 58477.5
$ws.Range("L83").Value = 118709.1
$ws.Range("M83").Value = -53797.5
$ws.Range("N83").Value = -128069.1
$ws.Range("H113").Value = 2601.6
$ws.Range("J113").Value = 2601.6
$ws.Range("L113").Value = 7804.799999999999
$ws.Range("N113").Value = -12144.8
$ws.Range("H136").Value = 2363.3333
$ws.Range("I136").Value = 1975
$ws.Range("K136").Value = 5925
$ws.Range("M136").Value = -825

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3824.0667
$ws.Range("I132").Value = 3618.4443
$ws.Range("J132").Value = 4132.5
$ws.Range("K132").Value = 10855.3329
$ws.Range("L132").Value = 12397.5
$ws.Range("M132").Value = -8325.332900000001
$ws.Range("N132").Value = -17457.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H32").Value = 566.6667
$ws.Range("I32").Value = 100
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 100
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = 217
$ws.Range("N32").Value = -2134
$ws.Range("H40").Value = 4850
$ws.Range("I40").Value = 4850
$ws.Range("K40").Value = 4850
$ws.Range("M40").Value = -4714
$ws.Range("H41").Value = 2500
$ws.Range("I41").Value = 2500
$ws.Range("K41").Value = 2500
$ws.Range("M41").Value = -2062
$ws.Range("H55").Value = 1178.5
$ws.Range("I55").Value = 1197.5
$ws.Range("J55").Value = 1174.7
$ws.Range("K55").Value = 1197.5
$ws.Range("L55").Value = 1174.7
$ws.Range("M55").Value = -1024.5
$ws.Range("N55").Value = -1520.7
$ws.Range("H132").Value = 5668
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1720.5834
$ws.Range("J81").Value = 2222
$ws.Range("L81").Value = 4444
$ws.Range("N81").Value = -6566
$ws.Range("H84").Value = 1720.5834
$ws.Range("J84").Value = 2222
$ws.Range("L84").Value = 22220
$ws.Range("N84").Value = -32828
$ws.Range("H96").Value = 1350
$ws.Range("I96").Value = 1200
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 173
$ws.Range("N96").Value = -4246
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N137").ClearContents()
